# 早晚香禮.xlsx - add PDF of pinyin/zhuyin
#
# The "row" sheet used to carry a redundant numeric helper column in B
# (a plain duplicate of column A, headed "rowZhuyin") while the text
# alignment keyword ("left"/"right"/"center"/"centerTitle") lived in E
# and the zhuyin phonetic guide lived in F.
#
# This edit retires the unused helper column: the alignment strings move
# left into B (replacing the numeric duplicate and its header), column E
# (the old alignment column) is deleted outright, and the zhuyin column
# shifts left to take its place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "align" column (E) left into B, overwriting the redundant
# numeric row-duplicate / "rowZhuyin" header that used to live there.
# Using Copy + PasteSpecial (values) instead of a literal Value assignment
# so the existing shared-string entries are reused instead of duplicated.
$ws.Range("E1:E74").Copy()
$ws.Range("B1:B74").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

# Drop the now-duplicated "align" column; this shifts the zhuyin column
# (old F) left into E.
$ws.Columns("E:E").Delete()

# Re-fit the columns that now hold different content.
$ws.Columns("A:C").AutoFit()
$ws.Columns("A:A").ColumnWidth = 3.6666666666666665
$ws.Columns("B:B").ColumnWidth = 9.833333333333334
$ws.Columns("C:C").ColumnWidth = 35.333333333333336

# Reset the view: scroll back to the top-left and select E65.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E65").Select()
